# Applies the "kleine Formataenderung Klassendiagramm+ Version mit Bildern" edit:
#  - nudges the "Klassendiagramm." and "Realisierungstools:" caption boxes
#  - nudges two small pictures (Grafik 63 / Grafik 69) near the bottom-right
#  - replaces the big class-diagram picture ("Grafik 26") with a slightly
#    resized copy of itself (same embedded image, new size/position),
#    appended at the end of the shape stack (as "Grafik 58")

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

# --- 1) "Klassendiagramm." label (Rechteck 22, id 23) ---
$klassendiagrammLabel = Get-ShapeById $shapes 23
$klassendiagrammLabel.Left = 14450127 / 12700.0
$klassendiagrammLabel.Top = 23909376 / 12700.0

# --- 2) "Realisierungstools:" label (Rechteck 23, id 24) ---
$realisierungstoolsLabel = Get-ShapeById $shapes 24
$realisierungstoolsLabel.Left = 14414066 / 12700.0
$realisierungstoolsLabel.Top = 24289897 / 12700.0

# --- 3) small picture "Grafik 63" (id 64) ---
$grafik63 = Get-ShapeById $shapes 64
$grafik63.Left = 18784463 / 12700.0
$grafik63.Top = 26606641 / 12700.0

# --- 4) small picture "Grafik 69" (id 70) ---
$grafik69 = Get-ShapeById $shapes 70
$grafik69.Left = 16213358 / 12700.0
$grafik69.Top = 26695704 / 12700.0

# --- 5) replace the class-diagram picture "Grafik 26" (id 27) with a
#        resized copy of itself, appended at the end of the shape list ---
$oldDiagramPic = Get-ShapeById $shapes 27

$dupRange = $oldDiagramPic.Duplicate()
$newDiagramPic = $dupRange.Item(1)
$newDiagramPic.Name = "Grafik 58"
$newDiagramPic.Left = 14449582 / 12700.0
$newDiagramPic.Top = 18937798 / 12700.0
$newDiagramPic.Width = 6319404 / 12700.0
$newDiagramPic.Height = 4947124 / 12700.0

$oldDiagramPic.Delete()
